# Apply updated cryptos list values (commit: 'Updated cryptos list on Fri Apr 26 04:48:35 UTC 2024 with GitHub Actions').
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to keep a Text number format so the value round-trips as a string,
# matching the source data (inline strings) in the workbook.
$textCells = @(
    "D5",
    "D6",
    "D9",
    "D11",
    "D14",
    "D20",
    "D21",
    "D22",
    "D24",
    "D25",
    "D27",
    "D28",
    "D29",
    "D31",
    "D33",
    "D34",
    "D36",
    "D37",
    "D39",
    "D40",
    "D41",
    "D45",
    "D48",
    "D51",
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$updates = [ordered]@{
    "D2" = '64.322.98'
    "E2" = '  -0.07%  '
    "D3" = '3.137.88'
    "E3" = '  -0.51%  '
    "E4" = '  -0.06%  '
    "D5" = '610.26'
    "E5" = '  +0.72%  '
    "D6" = '143.59'
    "E6" = '  -2.49%  '
    "E7" = '  -0.02%  '
    "D8" = '3.134.21'
    "E8" = '  -0.48%  '
    "D9" = '0.528'
    "E9" = '  +0.39%  '
    "E10" = '  +0.17%  '
    "D11" = '5.37'
    "E11" = '  -2.43%  '
    "E12" = '  +0.23%  '
    "E13" = '  +2.69%  '
    "D14" = '35.56'
    "E14" = '  -1.46%  '
    "D15" = '3.652.00'
    "E15" = '  -0.49%  '
    "E16" = '  +2.59%  '
    "D17" = '64.307.85'
    "E17" = '  -0.09%  '
    "D18" = '3.147.02'
    "E18" = '  -0.08%  '
    "E19" = '  -1.17%  '
    "D20" = '477.64'
    "E20" = '  -0.21%  '
    "D21" = '14.71'
    "E21" = '  +0.46%  '
    "D22" = '0.722'
    "E22" = '  +1.79%  '
    "E23" = '  +1.37%  '
    "D24" = '13.64'
    "E24" = '  -0.78%  '
    "D25" = '85.18'
    "E25" = '  +2.01%  '
    "E26" = '  -0.05%  '
    "B27" = 'RenderToken'
    "C27" = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    "D27" = '8.58'
    "E27" = '  +1.98%  '
    "B28" = 'PancakeSwap'
    "C28" = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    "D28" = '2.79'
    "E28" = '  -2.87%  '
    "D29" = '7.38'
    "E29" = '  +8.11%  '
    "E30" = '  +2.81%  '
    "D31" = '2.09'
    "E31" = '  -4.75%  '
    "E32" = '  -0.14%  '
    "D33" = '26.66'
    "E33" = '  +1.72%  '
    "D34" = '2.63'
    "E34" = '  -4.47%  '
    "E35" = '  +0.60%  '
    "D36" = '5.95'
    "E36" = '  -0.64%  '
    "D37" = '52.46'
    "E37" = '  -3.50%  '
    "E38" = '  +4.51%  '
    "D39" = '455.04'
    "E39" = '  +1.88%  '
    "D40" = '3.01'
    "E40" = '  +4.68%  '
    "D41" = '0.0396'
    "E41" = '  +0.11%  '
    "E42" = '  -0.36%  '
    "E43" = '  -1.24%  '
    "D44" = '2.862.72'
    "E44" = '  +0.99%  '
    "D45" = '0.266'
    "E45" = '  -0.50%  '
    "E46" = '  -0.07%  '
    "E47" = '  +4.69%  '
    "D48" = '26.45'
    "E48" = '  +0.12%  '
    "E50" = '  -0.13%  '
    "D51" = '120.41'
    "E51" = '  +2.05%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
